$d = $word.ActiveDocument

function Split-WithHiddenBookmark($range) {
    # Splitting a Range into a standalone run (even when its formatting is
    # identical to its neighbours) requires forcing a real run boundary.
    # Adding + immediately deleting a bookmark over the range achieves this
    # (see iron_docx::om_apply::split_run_at, invoked from m_bookmarks_add)
    # without leaving any stray <w:rPr/> artifact behind.
    $name = "tmpSplitMark"
    $d.Bookmarks.Add($name, $range)
    $d.Bookmarks($name).Delete()
}

# --- 1. "Example the code for the view." -> "Examine" + " the code for the view." ---
$rng = $d.Content
$rng.Find.Execute("Example the code for the view.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$r1 = $d.Range($start, $start + 7)
$r1.Text = "Examine"
Split-WithHiddenBookmark $r1

# --- 2. "Example the code for the controller." -> "Examine " + "the code for the controller." ---
$rng2 = $d.Content
$rng2.Find.Execute("Example the code for the controller.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $rng2.Start
$r2 = $d.Range($start2, $start2 + 8)
$r2.Text = "Examine "
Split-WithHiddenBookmark $r2

# --- 3. "Add a label and input box to the form for the director." ->
#        "Add a label and input box to the form for the " + bold "director" + "." ---
$rng3 = $d.Content
$rng3.Find.Execute("Add a label and input box to the form for the director.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start3 = $rng3.Start
$full3 = $rng3.Text
$dirIdx = $full3.IndexOf("director")
$rBold = $d.Range($start3 + $dirIdx, $start3 + $dirIdx + 8)
$rBold.Bold = 1

# --- 4. Move the "_GoBack" bookmark from the end of the document to right
#        after the second "Run the application" run. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$runRng = $d.Content
$runRng.Find.Execute("Run the application", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runRng.Find.Execute("Run the application", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runEnd = $runRng.End

# Adding a bookmark at a position that sits exactly on a paragraph's final
# text boundary (just before its paragraph mark) mis-places it at the
# document start. Work around this by inserting a throwaway character
# after the target point, anchoring the bookmark just before it, then
# removing the throwaway character again.
$tmp = $d.Range($runEnd, $runEnd)
$tmp.InsertAfter("X")
$bmRange = $d.Range($runEnd, $runEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($runEnd, $runEnd + 1).Delete()
